$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "addCoupon" test row first, so the "yes" shared string
# stays referenced (by B5) while B2:B4 switch to "no" below. This keeps
# the shared-strings table ordering/indices identical to the authored
# workbook (addCoupon, then no, appended at the end).
$ws.Range("A5").Value = "addCoupon"
$ws.Range("B5").Value = "yes"
$ws.Range("C5").Value = "chrome"
$ws.Range("D5").Value = "'"
$ws.Range("E5").Value = "'"

# Flip execute flag from yes -> no for the existing login test rows.
$ws.Range("B2").Value = "no"
$ws.Range("B3").Value = "no"
$ws.Range("B4").Value = "no"

$ws.Range("B4").Select()
